$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 51.03333333333333
$ws.Range("H2").Value = 13.09052777777778
$ws.Range("I2").Value = 154.9486388888889
$ws.Range("J2").Value = 5141.523499999999

$ws.Range("G3").Value = 81.55555555555556
$ws.Range("H3").Value = 21.97744444444444
$ws.Range("I3").Value = 160.3030833333333
$ws.Range("J3").Value = 9320.026638888889

$ws.Range("G4").Value = 68.37777777777778
$ws.Range("H4").Value = 19.79555555555556
$ws.Range("I4").Value = 199.1854166666667
$ws.Range("J4").Value = 9001.01661111111

$ws.Range("G5").Value = 66.99259259259259
$ws.Range("H5").Value = 18.12261111111111
$ws.Range("I5").Value = 160.959787037037
$ws.Range("J5").Value = 7400.116046296297

$ws.Range("G6").Value = 116.0861111111111
$ws.Range("H6").Value = 33.21086111111111
$ws.Range("I6").Value = 174.2913888888889
$ws.Range("J6").Value = 14626.94758333333

$ws.Range("G7").Value = 100.325
$ws.Range("H7").Value = 27.31283333333333
$ws.Range("I7").Value = 146.3938888888889
$ws.Range("J7").Value = 11037.77788888889

$ws.Range("G8").Value = 122.2416666666667
$ws.Range("H8").Value = 33.24593518518518
$ws.Range("I8").Value = 151.8361944444445
$ws.Range("J8").Value = 13564.09412037037

$ws.Range("G9").Value = 253.7324074074074
$ws.Range("H9").Value = 68.29277777777777
$ws.Range("I9").Value = 160.4442314814815
$ws.Range("J9").Value = 27521.70260185185

$ws.Range("G10").Value = 143.2564814814815
$ws.Range("H10").Value = 35.26078703703703
$ws.Range("I10").Value = 101.1248055555556
$ws.Range("J10").Value = 12349.21887037037

$ws.Range("G11").Value = 143.8490740740741
$ws.Range("H11").Value = 38.5919537037037
$ws.Range("I11").Value = 136.7964351851852
$ws.Range("J11").Value = 14187.8535462963

$ws.Range("G12").Value = 219.0037037037037
$ws.Range("H12").Value = 50.60387037037037
$ws.Range("I12").Value = 70.51983333333334
$ws.Range("J12").Value = 16014.40355555556

$ws.Range("G13").Value = 185.6916666666667
$ws.Range("H13").Value = 50.73413888888889
$ws.Range("I13").Value = 160.4697222222222
$ws.Range("J13").Value = 20328.45875

$ws.Range("G14").Value = 107.9111111111111
$ws.Range("H14").Value = 26.92652777777778
$ws.Range("I14").Value = 115.81075
$ws.Range("J14").Value = 9411.079861111111

$ws.Range("G15").Value = 101.6916666666667
$ws.Range("H15").Value = 23.9125
$ws.Range("I15").Value = 113.3460555555556
$ws.Range("J15").Value = 8791.10163888889

$ws.Range("G16").Value = 151.5138888888889
$ws.Range("H16").Value = 30.18658333333333
$ws.Range("I16").Value = 59.98061111111112
$ws.Range("J16").Value = 9313.030583333333

$ws.Range("G17").Value = 73.05000000000001
$ws.Range("H17").Value = 18.65216666666667
$ws.Range("I17").Value = 151.50675
$ws.Range("J17").Value = 7283.675527777778

$ws.Range("G18").Value = 123.6694444444445
$ws.Range("H18").Value = 34.96922222222222
$ws.Range("I18").Value = 143.492
$ws.Range("J18").Value = 14383.65822222222

$ws.Range("G19").Value = 120.1092592592593
$ws.Range("H19").Value = 29.81715740740741
$ws.Range("I19").Value = 120.5108240740741
$ws.Range("J19").Value = 10487.19562962963
